# mushroom_calculator.xlsx - scenario value updates + active-sheet/selection bookkeeping
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Oyster sheet: bump weekly fruiting-room/environment inputs and target yield.
# Downstream cells (B19 sum, B33:D36, B39:B41) are formulas and recalc on
# their own once the inputs below change.
# ---------------------------------------------------------------------------
$oyster = $wb.Worksheets.Item("Oyster")
$oyster.Range("B14").Value = 5
$oyster.Range("B15").Value = 3
$oyster.Range("B16").Value = 1
$oyster.Range("B17").Value = 0
$oyster.Range("B22").Value = 35
$oyster.Range("B18").Select()

# ---------------------------------------------------------------------------
# Lions Mane sheet: same kind of update to its own inputs.
# ---------------------------------------------------------------------------
$lionsMane = $wb.Worksheets.Item("Lions Mane")
$lionsMane.Range("B14").Value = 5
$lionsMane.Range("B15").Value = 3
$lionsMane.Range("B16").Value = 1
$lionsMane.Range("B17").Value = 0
$lionsMane.Range("B22").Value = 50
$lionsMane.Range("F20").Select()

# ---------------------------------------------------------------------------
# Summary and Scenario Planner are fully formula-driven - no direct writes,
# just move the live selection to match where the author left off.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("I22").Select()

$scenarioPlanner = $wb.Worksheets.Item("Scenario Planner")
$scenarioPlanner.Range("M30").Select()

# Summary becomes the active tab (was "Scenario Planner").
$summary.Activate()
